# Converts an "RRGGBB" hex color string into the packed integer that the
# Word object model expects for Font.Color / Shading.BackgroundPatternColor
# (value = R + G*256 + B*65536).
function RGBColor($hexstr) {
    $r = [Convert]::ToInt32($hexstr.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hexstr.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hexstr.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$wordMlNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$enDash = [char]0x2013

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: add two trailing spaces to the existing sentence, then
#    append "(This is a change – Version for branch alternate)" as three
#    separate dark-red (C00000) runs, matching the target markup exactly.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)

$trailSpacePoint = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$trailSpacePoint.InsertAfter("  ")

$startPos = $p1.Range.End - 1
$ins1 = $d.Range($startPos, $startPos)
$ins1.InsertAfter("(This is a change ${enDash} Ve")
$endPos1 = $p1.Range.End - 1
$run1 = $d.Range($startPos, $endPos1)
$run1.Font.Color = RGBColor "C00000"

$ins2 = $d.Range($endPos1, $endPos1)
$ins2.InsertAfter("rsion for branch alternate")
$endPos2 = $p1.Range.End - 1
$run2 = $d.Range($endPos1, $endPos2)
$run2.Font.Color = RGBColor "C00000"

$ins3 = $d.Range($endPos2, $endPos2)
$ins3.InsertAfter(")")
$endPos3 = $p1.Range.End - 1
$run3 = $d.Range($endPos2, $endPos3)
$run3.Font.Color = RGBColor "C00000"

# ---------------------------------------------------------------------------
# 2) The bare empty paragraph right after "It will be treated..." becomes a
#    new empty, formatted paragraph (Calibri/bold/#202122 text, #F9F9F9
#    shading) that precedes the existing "The Raven" paragraph.
# ---------------------------------------------------------------------------
$pEmpty = $d.Paragraphs.Item(3)
$pEmpty.Range.InsertXML(
    "<w:p xmlns:w='$wordMlNs'>" +
        "<w:pPr>" +
            "<w:shd w:val='clear' w:color='auto' w:fill='F9F9F9'/>" +
            "<w:rPr>" +
                "<w:rFonts w:ascii='Calibri' w:eastAsia='Times New Roman' w:hAnsi='Calibri' w:cs='Calibri'/>" +
                "<w:b/>" +
                "<w:bCs/>" +
                "<w:color w:val='202122'/>" +
            "</w:rPr>" +
        "</w:pPr>" +
    "</w:p>"
)

# ---------------------------------------------------------------------------
# 3) Drop the trailing "ank God almighty, we are free at last." paragraph,
#    leaving a single empty paragraph right before the sectPr, as in the
#    target document.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
$pLast.Range.InsertXML("<w:p xmlns:w='$wordMlNs'/>")

# ---------------------------------------------------------------------------
# 4) Prune the styles that are no longer referenced by any content once the
#    paragraph above is gone (Heading2/4 [+ their linked Char styles],
#    Hyperlink, apple-converted-space, audio-tool, subscribe,
#    subscribe-more-info, generic-title, podcast-tools__subscribe-links).
#    Styles must be removed from the highest original index down to the
#    lowest, otherwise the host's internal style table can be indexed out of
#    bounds after earlier removals shift everything down.
# ---------------------------------------------------------------------------
$stylesToRemove = @(
    "podcast-tools__subscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading 4 Char",
    "Heading 2 Char",
    "Hyperlink",
    "apple-converted-space",
    "heading 4",
    "heading 2"
)
foreach ($styleName in $stylesToRemove) {
    $style = $d.Styles.Item($styleName)
    $style.Delete()
}

Write-Host "Edits applied"
